# repull data, push all data, mean calculation
# Update the "dSF" column (F) values for several rows to reflect the
# repulled / recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F5").Value = -3
$ws.Range("F10").Value = -5
$ws.Range("F16").Value = -7
$ws.Range("F18").Value = -2
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = -2
$ws.Range("F28").Value = -2
